# add DTCO prediction and interpolation
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

    $ws.Range("B1695").Value = 163.9888458251953
    $ws.Range("B1696").Value = 165.5909423828125
    $ws.Range("B1697").Value = 165.7267150878906
    $ws.Range("B1698").Value = 165.5909423828125
    $ws.Range("B1699").Value = 163.2541046142578
    $ws.Range("B1700").Value = 165.5909423828125
    $ws.Range("B1701").Value = 165.8069458007812
    $ws.Range("B1702").Value = 164.2048492431641
    $ws.Range("B1703").Value = 163.9888458251953
    $ws.Range("B1704").Value = 161.1752014160156
    $ws.Range("B1705").Value = 161.1752014160156
    $ws.Range("B1706").Value = 161.6975860595703
    $ws.Range("B1707").Value = 161.3420257568359
    $ws.Range("B1708").Value = 161.3420257568359
    $ws.Range("B1709").Value = 160.8187561035156
    $ws.Range("B1710").Value = 161.1752014160156
    $ws.Range("B1711").Value = 162.4671936035156
    $ws.Range("B1712").Value = 164.2836151123047
    $ws.Range("B1713").Value = 164.2414398193359
    $ws.Range("B1714").Value = 155.8772735595703
    $ws.Range("B1715").Value = 153.5084228515625
    $ws.Range("B1716").Value = 153.7683868408203
    $ws.Range("B1717").Value = 161.1884613037109
    $ws.Range("B1718").Value = 160.2864227294922
    $ws.Range("B1719").Value = 163.0125122070312
    $ws.Range("B1720").Value = 162.7883605957031
    $ws.Range("B1721").Value = 164.9007415771484
    $ws.Range("B1722").Value = 162.8925323486328
    $ws.Range("B1723").Value = 163.1207580566406
    $ws.Range("B1724").Value = 162.8555450439453
    $ws.Range("B1725").Value = 162.8555450439453
    $ws.Range("B1726").Value = 163.6801452636719
    $ws.Range("B1727").Value = 163.6379699707031
    $ws.Range("B1728").Value = 160.8146820068359
    $ws.Range("B1729").Value = 159.4221954345703
    $ws.Range("B1730").Value = 164.3559417724609
    $ws.Range("B1731").Value = 165.9580383300781
    $ws.Range("B1732").Value = 165.9580383300781
    $ws.Range("B1733").Value = 165.8857116699219
    $ws.Range("B1734").Value = 163.1482849121094
    $ws.Range("B1735").Value = 162.6029663085938
    $ws.Range("B1736").Value = 161.795166015625
    $ws.Range("B1737").Value = 161.795166015625
    $ws.Range("B1738").Value = 162.3404846191406
